$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(8)
$sub = $tr.Characters($para.Start, $para.Length)
$sub.Text = "Tune a model (a bit) to get better performance. "
